$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-410) holds a "Förändrad" date value that moved from
# 2023-09-06 (serial 45175) to 2023-09-08 (serial 45177) for every row.
$ws.Range("C2:C410").Value = 45177
